# RPA datasets push 2023-12-05
# A new IPO listing record (삼성스팩9호, listed 2023-12-04) needs to be
# inserted as the new row 2 of the data table; every existing data row
# (previously rows 2-30) shifts down by one row (rows 3-31). No other
# cell content changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that hold plain "YYYY-MM-DD" text (상장일/청약일/납입일). Excel's
# default type-inference would otherwise silently convert these back into
# real date serials when re-written via .Value, so they need the
# apostrophe-literal treatment to stay stored as shared-string text,
# exactly like the source file.
$dateCols = @("A", "O", "P")
$allCols  = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q")

# Shift every existing data row down by one, bottom-up so nothing gets
# clobbered before it is read.
for ($r = 30; $r -ge 2; $r--) {
    $nr = $r + 1
    foreach ($col in $allCols) {
        $srcCell = $ws.Range($col + $r)
        $dstCell = $ws.Range($col + $nr)
        $val = $srcCell.Value2
        if ($dateCols -contains $col) {
            $dstCell.Value = "'" + $val
            $dstCell.ClearFormats()
        } else {
            $dstCell.Value2 = $val
        }
    }
}

# Populate the freshly-opened row 2 with the new IPO record.
$ws.Range("A2").Value = "'2023-12-04"
$ws.Range("A2").ClearFormats()
$ws.Range("B2").Value2 = "삼성스팩9호"
$ws.Range("C2").Value2 = "코스닥"
$ws.Range("D2").Value2 = 200
$ws.Range("E2").Value2 = "삼성"
$ws.Range("F2").Value2 = 200
$ws.Range("G2").Value2 = "-"
$ws.Range("H2").Value2 = "-"
$ws.Range("I2").Value2 = "-"
$ws.Range("J2").Value2 = "-"
$ws.Range("K2").Value2 = "대표"
$ws.Range("L2").Value2 = "-"
$ws.Range("M2").Value2 = 2000
$ws.Range("N2").Value2 = 100
$ws.Range("O2").Value = "'2023-11-23"
$ws.Range("O2").ClearFormats()
$ws.Range("P2").Value = "'2023-11-28"
$ws.Range("P2").ClearFormats()
$ws.Range("Q2").Value2 = 7500000
